$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Fix the infinite loop bug: the 3rd search-round comparison letter was
#    wrongly "s" - it must be "t" to match the word being searched ("asi"->t,y,... )
$ws.Range("K8").Value2 = "t"

# 2. New explanatory cells to the right of the existing table (row 5-7)
$ws.Range("O5").Value2 = "vertical"
$ws.Range("P5").Value2 = "ccaracteres + cccaracteres-1 + 2"
$ws.Range("P6").Value2 = "100 + 99 + 2 = 201"
$ws.Range("P7").Value2 = "o"
$ws.Range("Q7").Value2 = "l"
$ws.Range("R7").Value2 = "a"
$ws.Range("S7").Value2 = "`n"

# 3. New rows 15-16 with a small "filas/columnas" example
$ws.Range("K15").Value2 = "f"
$ws.Range("L15").Value2 = 100
$ws.Range("N15").Value2 = "filas"
$ws.Range("O15").Value2 = 2

$ws.Range("K16").Value2 = "c"
$ws.Range("L16").Value2 = 2
$ws.Range("N16").Value2 = "columnas"

# 4. Apply the built-in "Neutral" cell style (yellow fill / brown font) plus
#    the thin border & left alignment that the cells already had, to the
#    byte-count result cells.
$countRng = $ws.Range("L5:M5,L8:M8,L11:M11")
$countRng.Style = "Neutral"

$ws.Range("E5").Copy() | Out-Null
foreach ($area in $countRng.Areas) {
    $area.PasteSpecial(-4122) | Out-Null
}
$countRng.Interior.Color = 10284031
$countRng.Font.Color = 22428
$excel.CutCopyMode = $false

# 5. Update the view: zoom to 150%, keep B1 as top-left cell and move the
#    active selection from I5 to G5.
$ws.Range("G5").Select() | Out-Null
$excel.ActiveWindow.TopLeftCell = $ws.Range("B1")
$excel.ActiveWindow.Zoom = 150
